$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore the original (pre-rebuild) values for columns A, B, C, E, F, G, rows 1-20.
# Column D is untouched by the diff (stays 0 throughout).

$ws.Range("A1").Value = 3507.4447058673791
$ws.Range("B1").Value = 75.080990153743926
$ws.Range("C1").Value = 0
$ws.Range("E1").Value = 0.79377726754779943
$ws.Range("F1").Value = 75.362153570055028
$ws.Range("G1").Value = 0

$ws.Range("A2").Value = 7014.8894117347581
$ws.Range("B2").Value = 262.9518941120773
$ws.Range("C2").Value = 0
$ws.Range("E2").Value = 0.87684184836498535
$ws.Range("F2").Value = 267.33970034175036
$ws.Range("G2").Value = 0

$ws.Range("A3").Value = 10522.334117602137
$ws.Range("B3").Value = 445.63334544097199
$ws.Range("C3").Value = 0
$ws.Range("E3").Value = 0.48658930159237646
$ws.Range("F3").Value = 450.61878350654467
$ws.Range("G3").Value = 0

$ws.Range("A4").Value = 11443.368154799755
$ws.Range("B4").Value = 601.65145423150386
$ws.Range("C4").Value = 0
$ws.Range("E4").Value = 1.0877419839311691
$ws.Range("F4").Value = 600.11834635928119
$ws.Range("G4").Value = 0

$ws.Range("A5").Value = 12364.402191997373
$ws.Range("B5").Value = 674.91227013021012
$ws.Range("C5").Value = 0
$ws.Range("E5").Value = 0.86050393301221462
$ws.Range("F5").Value = 661.00530322193788
$ws.Range("G5").Value = 0

$ws.Range("A6").Value = 13285.436229194991
$ws.Range("B6").Value = 712.52744106198145
$ws.Range("C6").Value = 0
$ws.Range("E6").Value = 1.5619198401650232
$ws.Range("F6").Value = 707.52516515421269
$ws.Range("G6").Value = 0

$ws.Range("A7").Value = 14206.470266392609
$ws.Range("B7").Value = 897.02125475805963
$ws.Range("C7").Value = 0
$ws.Range("E7").Value = 7.096745785814532
$ws.Range("F7").Value = 900.929181060818
$ws.Range("G7").Value = 0

$ws.Range("A8").Value = 14988.874867478238
$ws.Range("B8").Value = 951.89500961385511
$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 20.020907514052976
$ws.Range("F8").Value = 956.89784371394228
$ws.Range("G8").Value = 0

$ws.Range("A9").Value = 15771.279468563867
$ws.Range("B9").Value = 1129.9058511477929
$ws.Range("C9").Value = 0
$ws.Range("E9").Value = 101.43040186868569
$ws.Range("F9").Value = 1126.932271939618
$ws.Range("G9").Value = 0

$ws.Range("A10").Value = 16553.684069649498
$ws.Range("B10").Value = 1322.5539588092247
$ws.Range("C10").Value = 0
$ws.Range("E10").Value = 534.24347356630085
$ws.Range("F10").Value = 1316.0718886194265
$ws.Range("G10").Value = 0

$ws.Range("A11").Value = 20061.128775516878
$ws.Range("B11").Value = 1281.5790802270885
$ws.Range("C11").Value = 0
$ws.Range("E11").Value = 1542.9500583964277
$ws.Range("F11").Value = 1283.0114968889977
$ws.Range("G11").Value = 0

$ws.Range("A12").Value = 23568.573481384257
$ws.Range("B12").Value = 1482.7869402388103
$ws.Range("C12").Value = 0
$ws.Range("E12").Value = 6604.602080539853
$ws.Range("F12").Value = 1472.0526183223976
$ws.Range("G12").Value = 0

$ws.Range("A13").Value = 27076.018187251637
$ws.Range("B13").Value = 1506.6267644400207
$ws.Range("C13").Value = 0
$ws.Range("E13").Value = 17971.52102783402
$ws.Range("F13").Value = 1504.7258099243829
$ws.Range("G13").Value = 0

$ws.Range("A14").Value = 27997.052224449257
$ws.Range("B14").Value = 1842.0201416644545
$ws.Range("C14").Value = 0
$ws.Range("E14").Value = 36001.395799211859
$ws.Range("F14").Value = 1833.5557299837803
$ws.Range("G14").Value = 0

$ws.Range("A15").Value = 28918.086261646877
$ws.Range("B15").Value = 2024.238479100826
$ws.Range("C15").Value = 0
$ws.Range("E15").Value = 31819.135305917069
$ws.Range("F15").Value = 2020.596964415608
$ws.Range("G15").Value = 0

$ws.Range("A16").Value = 29839.120298844497
$ws.Range("B16").Value = 1872.3477302185245
$ws.Range("C16").Value = 0
$ws.Range("E16").Value = 18159.527291583348
$ws.Range("F16").Value = 1870.3391932280665
$ws.Range("G16").Value = 0

$ws.Range("A17").Value = 30760.154336042116
$ws.Range("B17").Value = 2116.9183348607548
$ws.Range("C17").Value = 0
$ws.Range("E17").Value = 36001.315091821481
$ws.Range("F17").Value = 2105.0614048940702
$ws.Range("G17").Value = 0

$ws.Range("A18").Value = 31542.558937127746
$ws.Range("B18").Value = 2176.0138139479095
$ws.Range("C18").Value = 0
$ws.Range("E18").Value = 36001.481649925219
$ws.Range("F18").Value = 2179.7880050249728
$ws.Range("G18").Value = 0

$ws.Range("A19").Value = 32324.963538213375
$ws.Range("B19").Value = 2550.150814866437
$ws.Range("C19").Value = 0
$ws.Range("E19").Value = 18230.52610004591
$ws.Range("F19").Value = 2542.0283127772204
$ws.Range("G19").Value = 0

$ws.Range("A20").Value = 33107.368139299004
$ws.Range("B20").Value = 2663.7922169485646
$ws.Range("C20").Value = 0
$ws.Range("E20").Value = 36001.693861633743
$ws.Range("F20").Value = 2662.0831581446364
$ws.Range("G20").Value = 0

# Restore the recorded selection from the original file (A1:G20 instead of A1:AJ20).
[void]$ws.Range("A1:G20").Select()
